$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.176.53"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3
$ws.Range("D3").Value = "3.546.48"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.41"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.65"
$ws.Range("E6").Value = "  -2.16%  "

# Row 7
$ws.Range("D7").Value = "3.547.54"
$ws.Range("E7").Value = "  +0.82%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("E10").Value = "  -4.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.02"
$ws.Range("E11").Value = "  +2.76%  "

# Row 12
$ws.Range("E12").Value = "  -2.30%  "

# Row 13
$ws.Range("D13").Value = "4.145.72"
$ws.Range("E13").Value = "  +0.76%  "

# Row 14
$ws.Range("E14").Value = "  -3.41%  "

# Row 15
$ws.Range("E15").Value = "  -3.59%  "

# Row 16
$ws.Range("D16").Value = "3.541.03"
$ws.Range("E16").Value = "  +0.75%  "

# Row 17
$ws.Range("D17").Value = "66.323.30"

# Row 18
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.96"
$ws.Range("E19").Value = "  +2.47%  "

# Row 20
$ws.Range("E20").Value = "  -2.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.09"
$ws.Range("E21").Value = "  -1.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.61"
$ws.Range("E22").Value = "  -2.41%  "

# Row 23
$ws.Range("E23").Value = "  -0.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.92"
$ws.Range("E24").Value = "  -1.00%  "

# Row 25
$ws.Range("D25").Value = "3.686.84"
$ws.Range("E25").Value = "  +0.73%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("E27").Value = "  -0.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.35"
$ws.Range("E28").Value = "  -4.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  -3.66%  "

# Row 30
$ws.Range("E30").Value = "  -0.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("E32").Value = "  -5.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.160"
$ws.Range("E33").Value = "  -4.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.40"
$ws.Range("E34").Value = "  +0.01%  "

# Row 35
$ws.Range("D35").Value = "3.534.13"
$ws.Range("E35").Value = "  +0.57%  "

# Row 36
$ws.Range("E36").Value = "  -1.72%  "

# Row 37
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.88"
$ws.Range("E38").Value = "  -2.11%  "

# Row 39
$ws.Range("E39").Value = "  -4.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.55"
$ws.Range("E41").Value = "  -0.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0864"
$ws.Range("E42").Value = "  -3.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.18"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.895"
$ws.Range("E44").Value = "  -0.19%  "

# Row 45
$ws.Range("E45").Value = "  -9.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.36"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.21"
$ws.Range("E47").Value = "  -6.93%  "

# Row 48
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.22"
$ws.Range("E48").Value = "  -8.03%  "

# Row 49
$ws.Range("E49").Value = "  -1.37%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.21"
$ws.Range("E50").Value = "  -3.63%  "

# Row 51
$ws.Range("E51").Value = "  -3.93%  "
